# Update "last_edited_time" (column D) text for the rows that were
# re-synced from Notion in this batch.
#   - rows 2-14 now carry the newer "2024-08-03T03:17:00.000Z" stamp
#     (rows 13 & 14 moved out of the "...T18:25..." group into the
#     "...T18:24..." group, whose text itself also advanced).
#   - rows 15-22 carry the "2024-08-03T03:18:00.000Z" stamp.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 4).Value = "2024-08-03T03:17:00.000Z"
}
for ($r = 15; $r -le 22; $r++) {
    $ws.Cells.Item($r, 4).Value = "2024-08-03T03:18:00.000Z"
}

# Recomputed attendance numbers after the upstream report refresh:
#   S  = Tổng công tại LONG XUYÊN.number
#   V  = Nửa ngày.number
#   AF = Đầy đủ.number
#   AI = Tổng công tại CẦN THƠ.number
#   AM = Tổng công.formula.number
#   AP = Nghỉ có phép.number
$updates = @{
    2 = @{ AF = 2; AI = 2; AM = 2 }
    3 = @{ AF = 2; AI = 2; AM = 2 }
    4 = @{ AF = 2; AI = 2; AM = 2 }
    5 = @{ AF = 2; AI = 2; AM = 2 }
    8 = @{ AF = 2; AI = 2; AM = 2 }
    9 = @{ AF = 2; AI = 2; AM = 2 }
    10 = @{ S = 2; AF = 2; AM = 2 }
    11 = @{ AF = 2; AI = 2; AM = 2 }
    13 = @{ S = 1; AF = 3; AI = 2; AM = 3 }
    14 = @{ AF = 2; AI = 2; AM = 2 }
    16 = @{ S = 2; AF = 2; AM = 2 }
    18 = @{ V = 1; AI = 1.5; AM = 1.5; AP = 1 }
    19 = @{ S = 2; AF = 2; AM = 2 }
    20 = @{ AF = 2; AI = 2; AM = 2 }
    21 = @{ S = 2; AF = 2; AM = 2 }
    22 = @{ S = 2; AF = 2; AM = 2 }
}

foreach ($r in $updates.Keys) {
    foreach ($col in $updates[$r].Keys) {
        $ws.Range("$col$r").Value = $updates[$r][$col]
    }
}
